$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: Replace the first run's text ("(Ananth): Fig 3. ... multiplied")
# with the concatenation of the three new sentences that must precede it.
# The bookmark ("_GoBack") sits right after this run, so it naturally
# stays attached right after this newly-expanded text.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "(Ananth): Fig 3. Note that the signals are not added but rather multiplied",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Pulse Rate Frequency PRF (typical): 2 or 3 Mhz (500 or 333 nano sec interval)" +
    "Each pulse: 20 ns radio on; 15 ns radio off; 20 ns radio on; then off for rest of interval" +
    "At 2 MhZ, (20+15+20) / 500 = 11% ",
    2) | Out-Null

# ------------------------------------------------------------------
# Step 2: Split that big run into three separate paragraphs.
# These split points are far from the bookmark, so the bookmark is
# unaffected and simply stays glued to the end of the text that
# precedes it (the "At 2 MhZ..." sentence).
# ------------------------------------------------------------------
$split1 = $d.Content.Duplicate
$split1.Find.Execute("Each pulse:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p1 = $d.Range($split1.Start, $split1.Start)
$p1.InsertParagraphBefore() | Out-Null

$split2 = $d.Content.Duplicate
$split2.Find.Execute("At 2 MhZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p2 = $d.Range($split2.Start, $split2.Start)
$p2.InsertParagraphBefore() | Out-Null

# ------------------------------------------------------------------
# Step 3: Separate the remaining " so that -2 * -2 = +4." text (still
# attached right after the bookmark) into its own paragraph. Splitting
# exactly at the bookmark boundary would drag the bookmark along with
# it, so instead split one character later (inside the leading space)
# which leaves the bookmark correctly behind, then trim that now
# dangling leading character off of the "At 2 MhZ..." paragraph.
# ------------------------------------------------------------------
$split3 = $d.Content.Duplicate
$split3.Find.Execute(" so that -2 * -2 = +4.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $split3.Start + 1
$p3 = $d.Range($splitPoint, $splitPoint)
$p3.InsertParagraphBefore() | Out-Null

$mhzPara = $d.Content.Duplicate
$mhzPara.Find.Execute("At 2 MhZ, (20+15+20) / 500 = 11%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mhzParaRange = $d.Paragraphs.Item($mhzPara.Paragraphs.First.Index).Range
$trim = $d.Range($mhzParaRange.End - 2, $mhzParaRange.End - 1)
$trim.Text = ""

# ------------------------------------------------------------------
# Step 4: Insert a blank paragraph between the "At 2 MhZ..." paragraph
# and the trailing "so that..." text, matching the blank line in the
# target layout.
# ------------------------------------------------------------------
$soThat = $d.Content.Duplicate
$soThat.Find.Execute("so that -2 * -2 = +4.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$soThatParaIndex = $soThat.Paragraphs.First.Index
$soThatPara = $d.Paragraphs.Item($soThatParaIndex).Range
$soThatPara.InsertParagraphBefore() | Out-Null

# ------------------------------------------------------------------
# Step 5: Re-attach the "(Ananth): Fig 3. ..." lead-in text to the
# front of the trailing "so that..." fragment so it reads as a single
# sentence again, now in its own final paragraph.
# ------------------------------------------------------------------
$soThat2 = $d.Content.Duplicate
$soThat2.Find.Execute("so that -2 * -2 = +4.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$prependPoint = $d.Range($soThat2.Start, $soThat2.Start)
$prependPoint.InsertBefore("(Ananth): Fig 3. Note that the signals are not added but rather multiplied ") | Out-Null
